# Generate Report for Handoff
# Replace the old handoff GUID/filename ("e3056481-8f50-432b-86c6-10c020739433")
# and related hashed xlf names with the new ones, and bump the handoff
# timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldId = "e3056481-8f50-432b-86c6-10c020739433"
$newId = "7cb241d1-1a9b-4483-8afd-64024a79a70b"

$oldZhHash = "fea7338721c01c096ad60287103e87e47466093a"
$newZhHash = "bdacf8f57db88224d41e87ef38ac16b62fababf8"
$oldDeHash = "fea7338721c01c096ad60287103e87e47466093a"
$newDeHash = "bdacf8f57db88224d41e87ef38ac16b62fababf8"

$newMdName    = "$newId.md"
$newZhXlfName = "$newId.$newZhHash.zh-cn.xlf"
$newDeXlfName = "$newId.$newDeHash.de-de.xlf"

$newOverviewDate = "2016-17-17 11:17:18"
$newZhDate       = "2016-03-17 11:17:14"
$newDeDate       = "2016-03-17 11:17:18"

# ---------------------------------------------------------------------
# Sheet "Overview": A2 md file name, D2 handoff date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate

$ovA2 = $wsOverview.Range("A2")
$ovAddr = "https://github.com/OpenLocalizationTest/oltest/blob/2f50f7d019a6b28eebc8a55e0c86c0b5c91bd2ae/e2e/$oldId.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($ovA2, $ovAddr, "", "", $newMdName)

# ---------------------------------------------------------------------
# Sheet "zh-cn": A2 md file name, B2 (.md, unchanged link), D2 xlf file
# name, E2 handoff datetime
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName
$wsZh.Range("E2").Value = $newZhDate

$zhA2 = $wsZh.Range("A2")
$zhB2 = $wsZh.Range("B2")
$zhD2 = $wsZh.Range("D2")

$zhAddrMd  = "https://github.com/OpenLocalizationTest/oltest/blob/2f50f7d019a6b28eebc8a55e0c86c0b5c91bd2ae/e2e/$oldId.md"
$zhAddrXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/368dd85bcb4e963749f5c60538cb8aa246a74c84/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.$oldZhHash.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($zhA2, $zhAddrMd, "", "", $newMdName)
$wsZh.Hyperlinks.Add($zhB2, $zhAddrMd, "", "", ".md")
$wsZh.Hyperlinks.Add($zhD2, $zhAddrXlf, "", "", $newZhXlfName)

# ---------------------------------------------------------------------
# Sheet "de-de": A2 md file name, B2 (.md, unchanged link), D2 xlf file
# name, E2 handoff datetime
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName
$wsDe.Range("E2").Value = $newDeDate

$deA2 = $wsDe.Range("A2")
$deB2 = $wsDe.Range("B2")
$deD2 = $wsDe.Range("D2")

$deAddrMd  = "https://github.com/OpenLocalizationTest/oltest/blob/2f50f7d019a6b28eebc8a55e0c86c0b5c91bd2ae/e2e/$oldId.md"
$deAddrXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e717c8ed809d89ec07a79673d4737df7fd3cdf7f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.$oldDeHash.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($deA2, $deAddrMd, "", "", $newMdName)
$wsDe.Hyperlinks.Add($deB2, $deAddrMd, "", "", ".md")
$wsDe.Hyperlinks.Add($deD2, $deAddrXlf, "", "", $newDeXlfName)
